# Workbook: Jogos_do_Dia_Betfair_Back_Lay_2025-12-22
# This script applies betting-odds updates to Sheet1:
#  1) Updates several "Odd_*" values for the matches in rows 2-10 that keep their position.
#  2) Inserts a new row for a newly listed match (Dutch Eerste Divisie: Jong Ajax Amsterdam vs RKC Waalwijk)
#     at row 15, pushing the following matches down by one row.
#  3) Re-writes rows 11-20 completely with their final contents, which also accounts for the
#     Portuguese Segunda Liga / Romanian Liga I rows swapping order (rows 11 & 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update odds for matches that keep their row position (rows 2-10) ---
# Row 3: Turkish Super League | Basaksehir vs Gaziantep FK
$ws.Range("I3").Value = 5.2
$ws.Range("P3").Value = 2.2
$ws.Range("Q3").Value = 1.6

# Row 4: Azerbaijan Premier League | FK Sumqayit vs FC Sabah
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 1.84
$ws.Range("I4").Value = 2.06
$ws.Range("P4").Value = 1.51
$ws.Range("Q4").Value = 2.56

# Row 5: Greek Super League | Panserraikos vs Levadiakos
$ws.Range("H5").Value = 1.56
$ws.Range("I5").Value = 1.68
$ws.Range("P5").Value = 1.98
$ws.Range("Q5").Value = 1.84

# Row 6: Cypriot 1st Division | A.E.L. vs AEK Larnaca
$ws.Range("F6").Value = 3.35
$ws.Range("G6").Value = 5.9
$ws.Range("H6").Value = 1.8
$ws.Range("I6").Value = 2.06
$ws.Range("K6").Value = 4.2
$ws.Range("P6").Value = 1.88
$ws.Range("Q6").Value = 1.9

# Row 8: Turkish 1 Lig | Bodrum Belediyesi Bodru vs Amed Sportif Faaliyetle
$ws.Range("F8").Value = 1.96
$ws.Range("G8").Value = 2.34
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 4.8
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 4.2
$ws.Range("P8").Value = 2.02
$ws.Range("Q8").Value = 1.78

# Row 9: Turkish Super League | Genclerbirligi vs Trabzonspor
$ws.Range("F9").Value = 3.6
$ws.Range("K9").Value = 3.8
$ws.Range("P9").Value = 2.04
$ws.Range("Q9").Value = 1.78

# Row 10: Serbian Super League | FK Napredak vs FK Radnicki 1923
$ws.Range("F10").Value = 1.8
$ws.Range("H10").Value = 1.73
$ws.Range("I10").Value = 2.24
$ws.Range("J10").Value = 3.25
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 1.55

# --- Step 2: insert a new row at position 15 for the newly added match ---
# (this shifts the former rows 15-19 down to rows 16-20)
$ws.Range("A15").EntireRow.Insert()

# --- Step 3: rewrite rows 11-20 completely with their final values ---
# Row 11 becomes "Romanian Liga I" (previously at row 12) and row 12 becomes
# "Portuguese Segunda Liga" (previously at row 11); rows 13-14 keep their data (odds updated);
# row 15 is the brand-new match; rows 16,17,19,20 hold the former rows 15,16,17,19 (odds updated);
# row 18 (English Premier League) keeps its position (odds updated).

# Force the Date column (B) to be stored as text so "2025-12-22" is not auto-converted to a date
$ws.Range("B11:B20").NumberFormat = "@"

# Row 11: Romanian Liga I | Universitatea Craiova vs Csikszereda
$ws.Range("A11").Value = "Romanian Liga I"
$ws.Range("B11").Value = "2025-12-22"
$ws.Range("C11").Value = "15:00:00"
$ws.Range("D11").Value = "Universitatea Craiova"
$ws.Range("E11").Value = "Csikszereda"
$ws.Range("F11").Value = 1.04
$ws.Range("G11").Value = 1.4
$ws.Range("H11").Value = 1.04
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 5.4
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 2.28
$ws.Range("Q11").Value = 1.45
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AD11").Value = 0
$ws.Range("AE11").Value = 0
$ws.Range("AF11").Value = 0
$ws.Range("AG11").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AI11").Value = 0
$ws.Range("AJ11").Value = 0
$ws.Range("AK11").Value = 0
$ws.Range("AL11").Value = 0
$ws.Range("AM11").Value = 0
$ws.Range("AN11").Value = 0
$ws.Range("AO11").Value = 0

# Row 12: Portuguese Segunda Liga | Benfica B vs Sporting Lisbon B
$ws.Range("A12").Value = "Portuguese Segunda Liga"
$ws.Range("B12").Value = "2025-12-22"
$ws.Range("C12").Value = "15:00:00"
$ws.Range("D12").Value = "Benfica B"
$ws.Range("E12").Value = "Sporting Lisbon B"
$ws.Range("F12").Value = 2.34
$ws.Range("G12").Value = 2.6
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 3.4
$ws.Range("J12").Value = 3.35
$ws.Range("K12").Value = 3.85
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 1.93
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 0
$ws.Range("U12").Value = 0
$ws.Range("V12").Value = 0
$ws.Range("W12").Value = 0
$ws.Range("X12").Value = 0
$ws.Range("Y12").Value = 0
$ws.Range("Z12").Value = 0
$ws.Range("AA12").Value = 0
$ws.Range("AB12").Value = 0
$ws.Range("AC12").Value = 0
$ws.Range("AD12").Value = 0
$ws.Range("AE12").Value = 0
$ws.Range("AF12").Value = 0
$ws.Range("AG12").Value = 0
$ws.Range("AH12").Value = 0
$ws.Range("AI12").Value = 0
$ws.Range("AJ12").Value = 0
$ws.Range("AK12").Value = 0
$ws.Range("AL12").Value = 0
$ws.Range("AM12").Value = 0
$ws.Range("AN12").Value = 0
$ws.Range("AO12").Value = 0

# Row 13: Israeli Premier League | Maccabi Haifa vs Beitar Jerusalem
$ws.Range("A13").Value = "Israeli Premier League"
$ws.Range("B13").Value = "2025-12-22"
$ws.Range("C13").Value = "15:30:00"
$ws.Range("D13").Value = "Maccabi Haifa"
$ws.Range("E13").Value = "Beitar Jerusalem"
$ws.Range("F13").Value = 2.22
$ws.Range("G13").Value = 2.52
$ws.Range("H13").Value = 2.88
$ws.Range("I13").Value = 3.3
$ws.Range("J13").Value = 3.7
$ws.Range("K13").Value = 4.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 2.4
$ws.Range("Q13").Value = 1.57
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("W13").Value = 0
$ws.Range("X13").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("Z13").Value = 0
$ws.Range("AA13").Value = 0
$ws.Range("AB13").Value = 0
$ws.Range("AC13").Value = 0
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 0
$ws.Range("AF13").Value = 0
$ws.Range("AG13").Value = 0
$ws.Range("AH13").Value = 0
$ws.Range("AI13").Value = 0
$ws.Range("AJ13").Value = 0
$ws.Range("AK13").Value = 0
$ws.Range("AL13").Value = 0
$ws.Range("AM13").Value = 0
$ws.Range("AN13").Value = 0
$ws.Range("AO13").Value = 0

# Row 14: Portuguese Primeira Liga | Alverca vs Porto
$ws.Range("A14").Value = "Portuguese Primeira Liga"
$ws.Range("B14").Value = "2025-12-22"
$ws.Range("C14").Value = "15:45:00"
$ws.Range("D14").Value = "Alverca"
$ws.Range("E14").Value = "Porto"
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 14.5
$ws.Range("H14").Value = 1.36
$ws.Range("I14").Value = 1.39
$ws.Range("J14").Value = 4.9
$ws.Range("K14").Value = 5.6
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 1.9
$ws.Range("Q14").Value = 1.97
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 0
$ws.Range("X14").Value = 0
$ws.Range("Y14").Value = 0
$ws.Range("Z14").Value = 0
$ws.Range("AA14").Value = 0
$ws.Range("AB14").Value = 0
$ws.Range("AC14").Value = 0
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
$ws.Range("AF14").Value = 0
$ws.Range("AG14").Value = 0
$ws.Range("AH14").Value = 0
$ws.Range("AI14").Value = 0
$ws.Range("AJ14").Value = 0
$ws.Range("AK14").Value = 0
$ws.Range("AL14").Value = 0
$ws.Range("AM14").Value = 0
$ws.Range("AN14").Value = 0
$ws.Range("AO14").Value = 0

# Row 15: Dutch Eerste Divisie | Jong Ajax Amsterdam vs RKC Waalwijk
$ws.Range("A15").Value = "Dutch Eerste Divisie"
$ws.Range("B15").Value = "2025-12-22"
$ws.Range("C15").Value = "16:00:00"
$ws.Range("D15").Value = "Jong Ajax Amsterdam"
$ws.Range("E15").Value = "RKC Waalwijk"
$ws.Range("F15").Value = 3.75
$ws.Range("G15").Value = 6.4
$ws.Range("H15").Value = 1.61
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 3.75
$ws.Range("K15").Value = 8.8
$ws.Range("L15").Value = 1.01
$ws.Range("M15").Value = 1.01
$ws.Range("N15").Value = 2.36
$ws.Range("O15").Value = 1.15
$ws.Range("P15").Value = 2.36
$ws.Range("Q15").Value = 1.39
$ws.Range("R15").Value = 1.58
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 1.01
$ws.Range("U15").Value = 1.01
$ws.Range("V15").Value = 2
$ws.Range("W15").Value = 1.19
$ws.Range("X15").Value = 1000
$ws.Range("Y15").Value = 1000
$ws.Range("Z15").Value = 1000
$ws.Range("AA15").Value = 1000
$ws.Range("AB15").Value = 1000
$ws.Range("AC15").Value = 1000
$ws.Range("AD15").Value = 1000
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 1000
$ws.Range("AG15").Value = 1000
$ws.Range("AH15").Value = 1000
$ws.Range("AI15").Value = 1000
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 1000
$ws.Range("AL15").Value = 1000
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("AO15").Value = 1000

# Row 16: Dutch Eerste Divisie | Jong FC Utrecht vs Roda JC
$ws.Range("A16").Value = "Dutch Eerste Divisie"
$ws.Range("B16").Value = "2025-12-22"
$ws.Range("C16").Value = "16:00:00"
$ws.Range("D16").Value = "Jong FC Utrecht"
$ws.Range("E16").Value = "Roda JC"
$ws.Range("F16").Value = 2.64
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2.42
$ws.Range("I16").Value = 2.74
$ws.Range("J16").Value = 3.7
$ws.Range("K16").Value = 4.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 2.36
$ws.Range("Q16").Value = 1.6
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("T16").Value = 0
$ws.Range("U16").Value = 0
$ws.Range("V16").Value = 0
$ws.Range("W16").Value = 0
$ws.Range("X16").Value = 0
$ws.Range("Y16").Value = 0
$ws.Range("Z16").Value = 0
$ws.Range("AA16").Value = 0
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 0
$ws.Range("AD16").Value = 0
$ws.Range("AE16").Value = 0
$ws.Range("AF16").Value = 0
$ws.Range("AG16").Value = 0
$ws.Range("AH16").Value = 0
$ws.Range("AI16").Value = 0
$ws.Range("AJ16").Value = 0
$ws.Range("AK16").Value = 0
$ws.Range("AL16").Value = 0
$ws.Range("AM16").Value = 0
$ws.Range("AN16").Value = 0
$ws.Range("AO16").Value = 0

# Row 17: Italian Serie C | Union Brescia vs Inter Milan (Res)
$ws.Range("A17").Value = "Italian Serie C"
$ws.Range("B17").Value = "2025-12-22"
$ws.Range("C17").Value = "16:30:00"
$ws.Range("D17").Value = "Union Brescia"
$ws.Range("E17").Value = "Inter Milan (Res)"
$ws.Range("F17").Value = 1.61
$ws.Range("G17").Value = 1.98
$ws.Range("H17").Value = 2.02
$ws.Range("I17").Value = 15
$ws.Range("J17").Value = 3.25
$ws.Range("K17").Value = 950
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 1.57
$ws.Range("Q17").Value = 2.02
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("U17").Value = 0
$ws.Range("V17").Value = 0
$ws.Range("W17").Value = 0
$ws.Range("X17").Value = 0
$ws.Range("Y17").Value = 0
$ws.Range("Z17").Value = 0
$ws.Range("AA17").Value = 0
$ws.Range("AB17").Value = 0
$ws.Range("AC17").Value = 0
$ws.Range("AD17").Value = 0
$ws.Range("AE17").Value = 0
$ws.Range("AF17").Value = 0
$ws.Range("AG17").Value = 0
$ws.Range("AH17").Value = 0
$ws.Range("AI17").Value = 0
$ws.Range("AJ17").Value = 0
$ws.Range("AK17").Value = 0
$ws.Range("AL17").Value = 0
$ws.Range("AM17").Value = 0
$ws.Range("AN17").Value = 0
$ws.Range("AO17").Value = 0

# Row 18: English Premier League | Fulham vs Nottm Forest
$ws.Range("A18").Value = "English Premier League"
$ws.Range("B18").Value = "2025-12-22"
$ws.Range("C18").Value = "17:00:00"
$ws.Range("D18").Value = "Fulham"
$ws.Range("E18").Value = "Nottm Forest"
$ws.Range("F18").Value = 2.58
$ws.Range("G18").Value = 2.6
$ws.Range("H18").Value = 3.05
$ws.Range("I18").Value = 3.15
$ws.Range("J18").Value = 3.4
$ws.Range("K18").Value = 3.45
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 1.09
$ws.Range("N18").Value = 3.55
$ws.Range("O18").Value = 1.37
$ws.Range("P18").Value = 1.86
$ws.Range("Q18").Value = 2.14
$ws.Range("R18").Value = 1.32
$ws.Range("S18").Value = 3.9
$ws.Range("T18").Value = 1.85
$ws.Range("U18").Value = 2.1
$ws.Range("V18").Value = 0
$ws.Range("W18").Value = 0
$ws.Range("X18").Value = 12.5
$ws.Range("Y18").Value = 11.5
$ws.Range("Z18").Value = 21
$ws.Range("AA18").Value = 55
$ws.Range("AB18").Value = 10
$ws.Range("AC18").Value = 7.2
$ws.Range("AD18").Value = 13.5
$ws.Range("AE18").Value = 38
$ws.Range("AF18").Value = 16
$ws.Range("AG18").Value = 11.5
$ws.Range("AH18").Value = 18
$ws.Range("AI18").Value = 50
$ws.Range("AJ18").Value = 38
$ws.Range("AK18").Value = 29
$ws.Range("AL18").Value = 44
$ws.Range("AM18").Value = 1000
$ws.Range("AN18").Value = 26
$ws.Range("AO18").Value = 38

# Row 19: Spanish La Liga | Athletic Bilbao vs Espanyol
$ws.Range("A19").Value = "Spanish La Liga"
$ws.Range("B19").Value = "2025-12-22"
$ws.Range("C19").Value = "17:00:00"
$ws.Range("D19").Value = "Athletic Bilbao"
$ws.Range("E19").Value = "Espanyol"
$ws.Range("F19").Value = 1.84
$ws.Range("G19").Value = 1.86
$ws.Range("H19").Value = 5.2
$ws.Range("I19").Value = 5.5
$ws.Range("J19").Value = 3.65
$ws.Range("K19").Value = 3.75
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 1.48
$ws.Range("P19").Value = 1.68
$ws.Range("Q19").Value = 2.42
$ws.Range("R19").Value = 1.24
$ws.Range("S19").Value = 4.7
$ws.Range("T19").Value = 2.2
$ws.Range("U19").Value = 1.79
$ws.Range("V19").Value = 0
$ws.Range("W19").Value = 0
$ws.Range("X19").Value = 11
$ws.Range("Y19").Value = 14.5
$ws.Range("Z19").Value = 40
$ws.Range("AA19").Value = 190
$ws.Range("AB19").Value = 6.8
$ws.Range("AC19").Value = 8.4
$ws.Range("AD19").Value = 22
$ws.Range("AE19").Value = 1000
$ws.Range("AF19").Value = 9.8
$ws.Range("AG19").Value = 10.5
$ws.Range("AH19").Value = 26
$ws.Range("AI19").Value = 130
$ws.Range("AJ19").Value = 21
$ws.Range("AK19").Value = 23
$ws.Range("AL19").Value = 55
$ws.Range("AM19").Value = 230
$ws.Range("AN19").Value = 17.5
$ws.Range("AO19").Value = 180

# Row 20: Portuguese Primeira Liga | Benfica vs Famalicao
$ws.Range("A20").Value = "Portuguese Primeira Liga"
$ws.Range("B20").Value = "2025-12-22"
$ws.Range("C20").Value = "17:45:00"
$ws.Range("D20").Value = "Benfica"
$ws.Range("E20").Value = "Famalicao"
$ws.Range("F20").Value = 1.33
$ws.Range("G20").Value = 1.37
$ws.Range("H20").Value = 10.5
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 5.6
$ws.Range("K20").Value = 5.9
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 2.12
$ws.Range("Q20").Value = 1.78
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("W20").Value = 0
$ws.Range("X20").Value = 0
$ws.Range("Y20").Value = 0
$ws.Range("Z20").Value = 0
$ws.Range("AA20").Value = 0
$ws.Range("AB20").Value = 0
$ws.Range("AC20").Value = 0
$ws.Range("AD20").Value = 0
$ws.Range("AE20").Value = 0
$ws.Range("AF20").Value = 0
$ws.Range("AG20").Value = 0
$ws.Range("AH20").Value = 0
$ws.Range("AI20").Value = 0
$ws.Range("AJ20").Value = 0
$ws.Range("AK20").Value = 0
$ws.Range("AL20").Value = 0
$ws.Range("AM20").Value = 0
$ws.Range("AN20").Value = 0
$ws.Range("AO20").Value = 0
